$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.806.82'
$ws.Range('E2').Value = '  +4.14%  '
$ws.Range('D3').Value = '2.264.26'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.532'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.52%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.482'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.58'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.56%  '
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.59'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.73%  '
$ws.Range('D15').Value = '2.614.78'
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.19'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.47%  '
$ws.Range('D17').Value = '2.258.33'
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.761'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.44%  '
$ws.Range('D19').Value = '41.714.07'
$ws.Range('E19').Value = '  +4.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.26'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.95%  '
$ws.Range('D21').Value = '0.0₃0905'
$ws.Range('E21').Value = '  +1.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.74'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.73'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.16%  '
$ws.Range('E25').Value = '  +4.05%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +4.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.52%  '
$ws.Range('E29').Value = '  +11.21%  '
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0746'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.00%  '
$ws.Range('E37').Value = '  +1.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.78'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.54%  '
$ws.Range('E39').Value = '  +2.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.104'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.69%  '
$ws.Range('E41').Value = '  +2.00%  '
$ws.Range('E42').Value = '  +3.80%  '
$ws.Range('D43').Value = '2.062.34'
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.25'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0278'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.05'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.94%  '
$ws.Range('E49').Value = '  +3.69%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.86'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.43%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.74%  '
